# Remove the decorative clutter that accumulated across the doc:
#   - inline screenshot images (InlineShapes)
#   - "────…" separator-line paragraphs
#   - empty paragraphs used only for extra spacing (w:before=40 twips / 2pt)
#
# Strategy: walk the Paragraphs collection once, identify every paragraph
# that matches one of the three patterns, collect them, then delete from
# the end of the document backwards so earlier indices stay valid.

$d = $word.ActiveDocument

$targets = New-Object System.Collections.ArrayList

$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1

    $isSeparator = $false
    $txt = $p.Range.Text
    if ($txt.Length -gt 0 -and $txt.Substring(0, 1) -eq [char]0x2500) {
        $isSeparator = $true
    }

    $hasDrawing = $p.Range.InlineShapes.Count -gt 0

    $isSpacerEmpty = $false
    if ($txt -eq "`r") {
        if ($p.Format.SpaceBefore -eq 2) {
            $isSpacerEmpty = $true
        }
    }

    if ($isSeparator -or $hasDrawing -or $isSpacerEmpty) {
        [void]$targets.Add($i)
    }
}

# Delete highest index first so the remaining indices don't shift.
$sorted = $targets | Sort-Object -Descending
foreach ($idx in $sorted) {
    $d.Paragraphs($idx).Range.Delete()
}

Write-Output ("Removed paragraphs: " + $targets.Count)
